$d = $word.ActiveDocument

# Paragraph 1
$p1 = $d.Paragraphs(1).Range
$p1.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="44"/></w:rPr><w:t>Exploring the Marvelous World of Chemistry</w:t></w:r></w:p>')

# Paragraph 2
$p2 = $d.Paragraphs(2).Range
$p2.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Dr</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Alex Williams</w:t></w:r></w:p>')

# Paragraph 3
$p3 = $d.Paragraphs(3).Range
$p3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>alexwilliams@emailworld</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>com</w:t></w:r></w:p>')

# Paragraph 5
$p5 = $d.Paragraphs(5).Range
$p5.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Chemistry, the study of matter and its properties, is an essential part of our understanding of the natural world</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It is a fundamental science that has led to advancements in various fields, touching all aspects of our lives</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> In this essay, we delve into the captivating world of chemistry, unraveling its basic principles, exploring its diverse applications, and appreciating its profound impact on society</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>From the intricate interactions of atoms to the vastness of the universe, chemistry is an awe-inspiring tapestry that weaves together the fabric of our existence</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It enables us to understand the composition of matter, from the tiniest particles to the vast celestial bodies that grace our night sky</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Chemistry empowers us to manipulate and transform substances, leading to the development of materials, drugs, and technologies that drive progress and improve lives</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The study of chemistry provides a window into the inner workings of the world around us, offering explanations for phenomena that may seem mysterious or enigmatic</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> It unveils the secrets of chemical reactions, revealing the energetic interplay between atoms and molecules</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Chemistry empowers us to comprehend the complexities of life itself, shedding light on the biochemical processes that sustain all living organisms</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>')

# Paragraph 7
$p7 = $d.Paragraphs(7).Range
$p7.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>Chemistry is a captivating science that unlocks the secrets of matter and its properties</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> It provides a profound understanding of the composition of the world around us, the transformations it undergoes, and its profound impact on life</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Through chemistry, we uncover the intricate interactions of atoms, harness the power of chemical reactions, and develop innovative materials and technologies that drive progress and enhance our lives</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Chemistry stands as a testament to the boundless capacity of human knowledge and ingenuity in unraveling the mysteries of the universe</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>')

# Append new empty paragraph at end of document (paragraph 8)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endRange = $d.Range($lastPara.Range.End, $lastPara.Range.End)
$endRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

Write-Output $d.Paragraphs.Count